$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = -8
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = -2
$ws.Range("F12").Value = 4
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = -6
$ws.Range("F16").Value = -1
$ws.Range("F18").Value = 2
